$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "231×6=" "972×6="
Replace-Text "392×9=" "902×8="
Replace-Text "102×2=" "548×9="
Replace-Text "573×5=" "648×4="
Replace-Text "762×9=" "468×4="
Replace-Text "342×6=" "740×6="
Replace-Text "759×8=" "133×3="
Replace-Text "490×3=" "844×9="
Replace-Text "547×6=" "495×4="
Replace-Text "368×3=" "180×8="
Replace-Text "816×7=" "716×2="
Replace-Text "450×5=" "877×7="
Replace-Text "999×2=" "231×4="
Replace-Text "917×9=" "975×9="
Replace-Text "512×4=" "296×7="
Replace-Text "268×9=" "897×7="
Replace-Text "266×7=" "452×8="
Replace-Text "612×4=" "785×7="
Replace-Text "433×7=" "345×8="
Replace-Text "616×6=" "564×7="
Replace-Text "808×7=" "236×3="
Replace-Text "765×7=" "892×5="
Replace-Text "530×7=" "171×6="
Replace-Text "177×4=" "772×7="
Replace-Text "431×7=" "370×7="
